$d = $word.ActiveDocument

# Update the title / date line
$d.Content.Find.Execute("2025-12-20 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-12-21 Sunday", 2)

# Update the division problems in the table (Table 1), addressed by
# (row, column) so duplicate problem texts (e.g. "24÷2=") are each
# replaced with their own distinct new value rather than a global
# find/replace collapsing them.
$t = $d.Tables.Item(1)

$edits = @(
    @{ Row = 1;  Col = 1; New = "86÷7=" },
    @{ Row = 1;  Col = 2; New = "83÷6=" },
    @{ Row = 1;  Col = 3; New = "71÷5=" },
    @{ Row = 1;  Col = 4; New = "32÷5=" },
    @{ Row = 1;  Col = 5; New = "56÷4=" },

    @{ Row = 5;  Col = 1; New = "61÷9=" },
    @{ Row = 5;  Col = 2; New = "94÷8=" },
    @{ Row = 5;  Col = 3; New = "63÷9=" },
    @{ Row = 5;  Col = 4; New = "80÷5=" },
    @{ Row = 5;  Col = 5; New = "17÷3=" },

    @{ Row = 9;  Col = 1; New = "25÷7=" },
    @{ Row = 9;  Col = 2; New = "34÷5=" },
    @{ Row = 9;  Col = 3; New = "53÷7=" },
    @{ Row = 9;  Col = 4; New = "93÷4=" },
    @{ Row = 9;  Col = 5; New = "76÷6=" },

    @{ Row = 13; Col = 1; New = "75÷2=" },
    @{ Row = 13; Col = 2; New = "70÷9=" },
    @{ Row = 13; Col = 3; New = "26÷7=" },
    @{ Row = 13; Col = 4; New = "88÷9=" },
    @{ Row = 13; Col = 5; New = "99÷5=" },

    @{ Row = 17; Col = 1; New = "91÷3=" },
    @{ Row = 17; Col = 2; New = "40÷6=" },
    @{ Row = 17; Col = 3; New = "19÷2=" },
    @{ Row = 17; Col = 4; New = "89÷7=" },
    @{ Row = 17; Col = 5; New = "57÷4=" }
)

foreach ($e in $edits) {
    $cell = $t.Cell($e.Row, $e.Col)
    $cell.Range.Text = $e.New
}
